$d = $word.ActiveDocument
$bullet = [char]0x2022

# ------------------------------------------------------------------
# 1. Contact line: merge the split "github.com/<spellcheck>danielmartincraig</spellcheck>"
#    runs (and the trailing linkedin run) back into a single run that also
#    keeps the github handle alongside the linkedin handle.
# ------------------------------------------------------------------
$oldContact = "(803)389-6750 " + $bullet + " danielmartincraig@gmail.com " + $bullet + " github.com/danielmartincraig " + $bullet + " linkedin.com/danielcraig23"
$d.Content.Find.Execute($oldContact, $true, $false, $false, $false, $false, $true, 1, $false, $oldContact, 2) | Out-Null

# ------------------------------------------------------------------
# 2. Web Engineering bullet: merge the "Web " / "Engineering" / " I and II"
#    runs (and drop the gramStart/gramEnd proofErr markers) into one run.
# ------------------------------------------------------------------
$oldWebEng = "Web Engineering I and II"
$d.Content.Find.Execute($oldWebEng, $true, $false, $false, $false, $false, $true, 1, $false, $oldWebEng, 2) | Out-Null

# ------------------------------------------------------------------
# 3. Insert a new "OBJECTIVE:" heading paragraph right after the contact
#    line and before "EDUCATION:". The existing _GoBack bookmark (originally
#    sitting at the very end of the document, after "Fluent in Spanish")
#    moves here, matching real Word's "last edit" tracking.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$eduPar = $d.Paragraphs(4)
$eduPar.Range.InsertParagraphBefore()

$objPar = $d.Paragraphs(4)
$objPar.Style = "Heading1"
$objPar.Range.Text = "OBJECTIVE: "

$pos = $objPar.Range.End - 1
$seg1 = "Eager to drive back-end solutions at"
$d.Range($objPar.Range.Start, $pos).InsertAfter($seg1) | Out-Null
$d.Range($pos, $pos + $seg1.Length).Font.Size = 12
$pos = $pos + $seg1.Length

$seg2 = " Southwest"
$d.Range($objPar.Range.Start, $pos).InsertAfter($seg2) | Out-Null
$d.Range($pos, $pos + $seg2.Length).Font.Size = 12
$pos = $pos + $seg2.Length

$seg3 = " on a full-time basis"
$d.Range($objPar.Range.Start, $pos).InsertAfter($seg3) | Out-Null
$d.Range($pos, $pos + $seg3.Length).Font.Size = 12

$d.Bookmarks.Add("_GoBack", $d.Range($pos, $pos)) | Out-Null
